# Apply the "counterparty_notes" sheet addition + counterparties eik_egn
# text-format fix, as described by the commit "Forgot to add some files".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Fix the `eik_egn` column on the "counterparties" sheet: it used to
#    store values like `"111222333"` (quoted literal text) so Excel would
#    not coerce them to numbers. Replace that hack with a real text
#    number format ("@") and store the plain digit string instead.
# ---------------------------------------------------------------------
$wsCounterparties = $wb.Worksheets.Item("counterparties")

$wsCounterparties.Columns.Item(2).NumberFormat = "@"

$eikValues = @("111222333", "111222334", "111222335", "111222336", "111222337", "111222338", "111222339")
for ($i = 0; $i -lt $eikValues.Length; $i++) {
    $row = $i + 2
    $wsCounterparties.Cells.Item($row, 2).Value = $eikValues[$i]
}

# ---------------------------------------------------------------------
# 2. Add the missing "counterparty_notes" worksheet after "counterparties"
#    with its header row and sample note rows.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNotes = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNotes.Name = "counterparty_notes"

$wsNotes.Range("A1").Value = "counterparty_id"
$wsNotes.Range("B1").Value = "note"
$wsNotes.Range("C1").Value = "created_by"
$wsNotes.Range("D1").Value = "updated_by"

$noteRows = @(
    @(1, "Note 1"),
    @(1, "Note 2"),
    @(1, "Note 3"),
    @(2, "Note 2 – 1"),
    @(2, "Note 2 – 2"),
    @(1, "Note 2 – 3"),
    @(3, "Note 3 – 1")
)

for ($i = 0; $i -lt $noteRows.Length; $i++) {
    $row = $i + 2
    $wsNotes.Cells.Item($row, 1).Value = $noteRows[$i][0]
    $wsNotes.Cells.Item($row, 2).Value = $noteRows[$i][1]
    $wsNotes.Cells.Item($row, 3).Value = 1
    $wsNotes.Cells.Item($row, 4).Value = 1
}

# ---------------------------------------------------------------------
# 3. Restore sane single-cell selections on every sheet (the previous
#    save had left multi-range leftover selections like "C5:D8 B6") and
#    make the new sheet the active tab, matching the final workbook view.
# ---------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("users")
$wsUsers.Activate() | Out-Null
$wsUsers.Range("C8").Select() | Out-Null

$wsRoles = $wb.Worksheets.Item("roles")
$wsRoles.Activate() | Out-Null
$wsRoles.Range("B6").Select() | Out-Null

$wsRolesUsers = $wb.Worksheets.Item("roles_users")
$wsRolesUsers.Activate() | Out-Null
$wsRolesUsers.Range("N17").Select() | Out-Null

$wsCompanyUnits = $wb.Worksheets.Item("company_units")
$wsCompanyUnits.Activate() | Out-Null
$wsCompanyUnits.Range("D33").Select() | Out-Null

$wsEmployees = $wb.Worksheets.Item("employees")
$wsEmployees.Activate() | Out-Null
$wsEmployees.Range("C24").Select() | Out-Null

$wsCounterparties.Activate() | Out-Null
$wsCounterparties.Range("D1").Select() | Out-Null

$wsNotes.Activate() | Out-Null
$wsNotes.Range("B9").Select() | Out-Null
